function Replace-ParagraphXml($Doc, $Anchor, $NewParagraphXml) {
    $r = $Doc.Content
    $found = $r.Find.Execute($Anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor not found: $Anchor"
    }
    # Clear the matched text, leaving a collapsed insertion point where the
    # replacement paragraph content should land.
    $r.Text = ""
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
           $NewParagraphXml +
           '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg)
}

$d = $word.ActiveDocument

# "House minimum height" -> "Building" + " minimum height" (two runs, matching
# the authored edit where "House" was swapped for "Building" in its own run).
Replace-ParagraphXml $d "House minimum height" (
    '<w:p><w:pPr><w:pStyle w:val="Liststycke"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Building</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> minimum height</w:t></w:r>' +
    '</w:p>'
)

# "House maximum height" -> "Building" + " maximum height", and this is also
# where the stray "_GoBack" bookmark ends up (moved from the paragraph below).
Replace-ParagraphXml $d "House maximum height" (
    '<w:p><w:pPr><w:pStyle w:val="Liststycke"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Building</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> maximum height</w:t></w:r>' +
    '</w:p>'
)

# Drop the old "_GoBack" bookmark from the end of the district-variables
# paragraph (it moved to the "Building maximum height" bullet above).
Replace-ParagraphXml $d "The three last variables are connected to specific districts. This means that three districts, with three variables each, the user can control nine different variables controlling the city generation plus the seed. The seed is a string of characters. The seed will change how the city is generated. This is needed so the user can generate several cities with the same district parameters but different seeds. The seed also makes the generation deterministic meaning that the exact same city can be generated if all the parameters are entered. " (
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr>' +
    '<w:t xml:space="preserve">The three last variables are connected to specific districts. This means that three districts, with three variables each, the user can control nine different variables controlling the city generation plus the seed. The seed is a string of characters. The seed will change how the city is generated. This is needed so the user can generate several cities with the same district parameters but different seeds. The seed also makes the generation deterministic meaning that the exact same city can be generated if all the parameters are entered. </w:t>' +
    '</w:r></w:p>'
)

# Fix the "i.e" spelling-error markup: drop the proofErr wrapper and add the
# missing period so it reads "i.e.".
Replace-ParagraphXml $d "The city edges will not be realistic (i.e no smaller roads or villages at the edge of the city)" (
    '<w:p><w:pPr><w:pStyle w:val="Liststycke"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr>' +
    '<w:rPr><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>The city edges will not be realistic (</w:t></w:r>' +
    '<w:r><w:rPr><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>i.e.</w:t></w:r>' +
    '<w:r><w:rPr><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> no smaller roads or villages at the edge of the city)</w:t></w:r>' +
    '</w:p>'
)

Write-Output "Edits applied"
